# Commit: "Added three commmodites i.e w+frk, frk rra, frk br"
#
# 1. Remove the sample data rows from the "wheat", "rra" and "coarse_grain"
#    sheets, leaving only the header row on each.
# 2. Add three new sheets ("frk_rra", "frk_br", "frk") at the end of the
#    workbook, each with just the same header row used on every other
#    sheet (From / From State / To / To State / Commodity / Values).

$wb = $excel.ActiveWorkbook

# --- 1. Strip the data rows, keep the header row only -----------------
$ws = $wb.Worksheets.Item("wheat")
$ws.Rows("2:3").Delete()

$ws = $wb.Worksheets.Item("rra")
$ws.Rows("2:2").Delete()

$ws = $wb.Worksheets.Item("coarse_grain")
$ws.Rows("2:2").Delete()

# --- 2. Add the new commodity sheets -----------------------------------
$newSheetNames = @("frk_rra", "frk_br", "frk")

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

foreach ($name in $newSheetNames) {
    $newWs = $wb.Worksheets.Add($null, $lastSheet)
    $newWs.Name = $name

    $newWs.Range("A1").Value = "From"
    $newWs.Range("B1").Value = "From State"
    $newWs.Range("C1").Value = "To"
    $newWs.Range("D1").Value = "To State"
    $newWs.Range("E1").Value = "Commodity"
    $newWs.Range("F1").Value = "Values"

    $headerRange = $newWs.Range("A1:F1")
    $headerRange.Font.Bold = $true
    $headerRange.Borders.LineStyle = 1
    $headerRange.HorizontalAlignment = -4108
    $headerRange.VerticalAlignment = -4160

    $lastSheet = $newWs
}
